# Configurable CV voltage excitation sequence
# Adds a new "Set slope duration between Voltage1 & Voltage2 for CV signal [us]"
# sub-command (ID 103) to the "GUI to MCU" sheet, right before the existing
# "Exit parameter configuration" row, shifting everything below it down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GUI to MCU")

# 1) Insert a new row at position 11 (pushes old row 11 "Exit parameter
#    configuration" down to row 12, and the drawing shapes below it too).
$ws.Rows.Item(11).Insert()

# 2) Copy the cell formatting from row 10 (the row directly above) into the
#    F:H cells of the new row 11 so the styles line up with the other
#    "sub command" rows (value/returns/description columns).
$ws.Range("F10:H10").Copy()
$ws.Range("F11:H11").PasteSpecial(-4122)

# Match the row height used by the sibling sub-command rows (8, 9, 10, 12).
$ws.Rows.Item(11).RowHeight = 30

# 3) Fill in the new command's data.
$ws.Range("C11").Value = "Set slope duration between Voltage1 & Voltage2 for CV signal [us]"
$ws.Range("E11").Value = 103
$ws.Range("F11").Value = "value (4 bytes)"
$ws.Range("G11").Value = "success bool (1 byte):`n0 if cmd ID not recognized"
$ws.Range("H11").Value = "Set slope duration between voltage 1 & 2 for the CV signal"

# 4) The drawing shapes anchored below the insertion point don't automatically
#    move with this engine's Rows.Insert(), so reposition them explicitly to
#    their new (shifted-by-one-row) cell anchors, preserving their sizes.
$arrow = $ws.Shapes.Item(1)
$arrow.Top = 368.25
$arrow.Height = 104.25

$textbox = $ws.Shapes.Item(2)
$textbox.Top = 481.5000787401575
$textbox.Height = 242.71149606299207
